$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45178 -> 45179) for every data row (rows 2 through 307).
$ws.Range("C2:C307").Value = 45179
